# Zerar lista de musicas: atualiza as 2 musicas existentes e adiciona
# mais 6 slides (2 musicas completas de 4 slides cada).

$p = $ppt.ActivePresentation

# --- Slide 1 : "DEUS TRINO DE AMOR" -> "GLORIA BANDA CIROS" (estrofe 1) ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "GLORIA BANDA CIROS"
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "Estrofe 1: Deus e Pai nós vos louvamos`rAdoramos, bendizemos`rDamos glória ao vosso nome`rVossos dons agradecemos"

# --- Slide 2 : Refrão da nova música ---
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Refrão: Glória, glória a Deus no céu e paz na Terra"

# --- Slide 3 (novo) : Estrofe 2 - duplicado a partir do slide 2 (corpo unico) ---
$last = $p.Slides.Item($p.Slides.Count)
$s3 = $last.Duplicate()
$s3.Name = "Slide 3"
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Estrofe 2: Senhor nosso Jesus Cristo`rUnigênito do Pai`rVoz de Deus, Cordeiro Santo`rNossas culpas perdoai"

# --- Slide 4 (novo) : Refrão repetido ---
$last = $p.Slides.Item($p.Slides.Count)
$s4 = $last.Duplicate()
$s4.Name = "Slide 4"
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Refrão: Glória, glória a Deus no céu e paz na Terra"

# --- Slide 5 (novo) : nova música "UM CORAÇÃO PARA AMAR" - duplicado a partir do slide 1 (titulo + corpo) ---
$last = $p.Slides.Item($p.Slides.Count)
$s5 = $last.Duplicate()
$s5.Name = "Slide 5"
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "UM CORAÇÃO PARA AMAR"
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "Estrofe 1: Um coração para amar, pra perdoar e sentir`rPara chorar e sorrir, ao me criar Tu me destes`rUm coração pra sonhar, inquieto e sempre a bater`rAnsioso por entender as coisas que Tu disseste"

# --- Slide 6 (novo) : Refrão - duplicado a partir do slide 2 (corpo unico) ---
$last = $p.Slides.Item($p.Slides.Count)
$s6 = $last.Duplicate()
$s6.Name = "Slide 6"
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Refrão: Eis o que eu venho Te dar`rEis o que eu ponho no altar`rToma, Senhor, que ele é Teu`rMeu coração não é meu"

# --- Slide 7 (novo) : Estrofe 2 ---
$last = $p.Slides.Item($p.Slides.Count)
$s7 = $last.Duplicate()
$s7.Name = "Slide 7"
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Estrofe 2: Quero que o meu coração seja tão cheio de paz`rQue não se sinta capaz de sentir ódio ou rancor`rQuero que a minha oração possa me amadurecer`rLeve-me a compreender as consequências do amor"

# --- Slide 8 (novo) : Refrão repetido ---
$last = $p.Slides.Item($p.Slides.Count)
$s8 = $last.Duplicate()
$s8.Name = "Slide 8"
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Refrão: Eis o que eu venho Te dar`rEis o que eu ponho no altar`rToma, Senhor, que ele é Teu`rMeu coração não é meu"

Write-Output ("Slides count: " + $p.Slides.Count)
